$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: D8 and E8
$ws.Range("D8").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E8").Value = "['Normal', 'HardwareFault']"

# Row 59: D59 and E59
$ws.Range("D59").Value = "[0, 1, 0, 0, 0, 0, 0]"
$ws.Range("E59").Value = "['SurroundingEnvironment']"

# Row 88: D88 and E88
$ws.Range("D88").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal']"

# Row 92: D92 and E92
$ws.Range("D92").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E92").Value = "['Normal', 'SoftwareFault']"

# Row 113: D113 and E113
$ws.Range("D113").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal']"
